$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.164.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.789.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("E6").Value = "  +2.10%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.74"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.23"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("E10").Value = "  +1.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0665"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.045.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +13.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.777.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.07%  "
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "34.137.61"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "255.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0746"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.98%  "
$ws.Range("E25").Value = "  -1.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("E31").Value = "  +0.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0518"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.36%  "
$ws.Range("E34").Value = "  +2.60%  "
$ws.Range("E35").Value = "  +2.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.454.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.04%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  +3.59%  "
$ws.Range("E39").Value = "  +0.87%  "
$ws.Range("E40").Value = "  +1.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "83.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0510"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.81%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("E47").Value = "  +3.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.944.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.98%  "
$ws.Range("E49").Value = "  +8.34%  "
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.18%  "
